# Auto-generated edit script: updates cryptos price/volume table
# to match the Sat Jan 20 17:49:57 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.519.32"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.457.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.05%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.74%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.09"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("E7").Value = "  +2.56%  "

$ws.Range("E8").Value = "  -0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.507"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.57%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0796"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.52%  "

$ws.Range("E12").Value = "  +1.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.836.80"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.65%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.457.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.767"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.38%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.505.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0934"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.15%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.58%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.08%  "

$ws.Range("E28").Value = "  +1.68%  "

$ws.Range("E29").Value = "  +1.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.64%  "

$ws.Range("E32").Value = "  +2.48%  "

$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.27%  "

$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.29%  "

$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.967.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0280"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.16%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.695.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.36%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.171"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.96%  "
